$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 12 new qubit-result rows (519-530) below the existing data (which
# ends at row 518), mirroring the "2019-11-15" qubit run.
# ---------------------------------------------------------------------------

# 1) Seed formatting for the new rows by copying the last existing data row
#    (518) down. Copy A:N and R separately so we don't introduce the blank
#    O/P/Q cells that a full A:R copy would create (row 518 has no values in
#    those columns).
$ws.Range("A518:N518").Copy($ws.Range("A519:N530"))
$ws.Range("R518").Copy($ws.Range("R519:R530"))

# 2) Column I in these new rows uses a different cell style (s="1", the
#    plain Calibri-12 style) than the rest of the legacy rows (s="2", Arial
#    10). Re-stamp the number format by copying from an existing cell that
#    already carries that style, without touching the values we set below.
$ws.Range("A263").Copy($ws.Range("I519:I530"))

# 3) Per-row data: run_ID, tube label, test_date, qubit conc, orig conc,
#    dilution reading (column I).
$ws.Cells.Item(519,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(519,3).Value2 = "Sample_#191115-145728"
$ws.Cells.Item(519,4).Value2 = 43784.623240740744
$ws.Cells.Item(519,5).Value2 = 157
$ws.Cells.Item(519,6).Value2 = 15.7
$ws.Cells.Item(519,9).Value2 = 91

$ws.Cells.Item(520,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(520,3).Value2 = "Sample_#191115-145719"
$ws.Cells.Item(520,4).Value2 = 43784.623136574075
$ws.Cells.Item(520,5).Value2 = 312
$ws.Cells.Item(520,6).Value2 = 31.2
$ws.Cells.Item(520,9).Value2 = 53

$ws.Cells.Item(521,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(521,3).Value2 = "Sample_#191115-145708"
$ws.Cells.Item(521,4).Value2 = 43784.62300925926
$ws.Cells.Item(521,5).Value2 = 228
$ws.Cells.Item(521,6).Value2 = 22.8
$ws.Cells.Item(521,9).Value2 = 124

$ws.Cells.Item(522,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(522,3).Value2 = "Sample_#191115-145700"
$ws.Cells.Item(522,4).Value2 = 43784.622916666667
$ws.Cells.Item(522,5).Value2 = 399
$ws.Cells.Item(522,6).Value2 = 39.9
$ws.Cells.Item(522,9).Value2 = 111

$ws.Cells.Item(523,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(523,3).Value2 = "Sample_#191115-145649"
$ws.Cells.Item(523,4).Value2 = 43784.622789351852
$ws.Cells.Item(523,5).Value2 = 282
$ws.Cells.Item(523,6).Value2 = 28.2
$ws.Cells.Item(523,9).Value2 = 24

$ws.Cells.Item(524,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(524,3).Value2 = "Sample_#191115-145640"
$ws.Cells.Item(524,4).Value2 = 43784.622685185182
$ws.Cells.Item(524,5).Value2 = 630
$ws.Cells.Item(524,6).Value2 = 63
$ws.Cells.Item(524,9).Value2 = 177

$ws.Cells.Item(525,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(525,3).Value2 = "Sample_#191115-145630"
$ws.Cells.Item(525,4).Value2 = 43784.622569444444
$ws.Cells.Item(525,5).Value2 = 323
$ws.Cells.Item(525,6).Value2 = 32.3
$ws.Cells.Item(525,9).Value2 = 119

$ws.Cells.Item(526,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(526,3).Value2 = "Sample_#191115-145622"
$ws.Cells.Item(526,4).Value2 = 43784.622476851851
$ws.Cells.Item(526,5).Value2 = 331
$ws.Cells.Item(526,6).Value2 = 33.1
$ws.Cells.Item(526,9).Value2 = 29

$ws.Cells.Item(527,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(527,3).Value2 = "Sample_#191115-145613"
$ws.Cells.Item(527,4).Value2 = 43784.622372685182
$ws.Cells.Item(527,5).Value2 = 229
$ws.Cells.Item(527,6).Value2 = 22.9
$ws.Cells.Item(527,9).Value2 = 10

$ws.Cells.Item(528,1).Value2 = "2019-11-15_145544"
$ws.Cells.Item(528,3).Value2 = "Sample_#191115-145604"
$ws.Cells.Item(528,4).Value2 = 43784.62226851852
$ws.Cells.Item(528,5).Value2 = 281
$ws.Cells.Item(528,6).Value2 = 28.1
$ws.Cells.Item(528,9).Value2 = 169

$ws.Cells.Item(529,1).Value2 = "2019-11-15_150415"
$ws.Cells.Item(529,3).Value2 = "Sample_#191115-150444"
$ws.Cells.Item(529,4).Value2 = 43784.628287037034
$ws.Cells.Item(529,5).Value2 = 221
$ws.Cells.Item(529,6).Value2 = 22.1
$ws.Cells.Item(529,9).Value2 = 28

$ws.Cells.Item(530,1).Value2 = "2019-11-15_150415"
$ws.Cells.Item(530,3).Value2 = "Sample_#191115-150435"
$ws.Cells.Item(530,4).Value2 = 43784.628182870372
$ws.Cells.Item(530,5).Value2 = 510
$ws.Cells.Item(530,6).Value2 = 51
$ws.Cells.Item(530,9).Value2 = 81

# All 12 new rows share the same source URL in column R.
for ($r = 519; $r -le 530; $r++) {
    $ws.Cells.Item($r,18).Value2 = "https://grace-ac.github.io/rna-qubit-day9-setof12/"
}

# 4) Column M keeps the "(F)*(L-G)" calc, extended down through the new rows.
for ($r = 519; $r -le 530; $r++) {
    $ws.Cells.Item($r,13).Formula = "=(F$r)*(L$r-G$r)"
}

# ---------------------------------------------------------------------------
# Reflect the scrolled viewport / new selection from the edit session.
# ---------------------------------------------------------------------------
[void]$ws.Range("Q526").Select()
